$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

$ws2.Range("A1").Value = 1
$ws2.Range("A2").Value = 2
$ws2.Range("A3").Value = 3
$ws2.Range("A4").Value = 4
$ws2.Range("A5").Value = 5

$ws2.Range("A7").Formula = "=SUM(A1:A5)"
$ws2.Range("A8").Formula = "=AVERAGE(A1:A5)"
$ws2.Range("A9").Formula = "=MAX(A1:A5)"
$ws2.Range("A10").Formula = "=MIN(A1:A5)"
$ws2.Range("A11").Formula = "=SUBTOTAL(6,A1:A5)"

$ws2.Rows.Item(11).Select()
